# Weekly price update: insert 4 new "Femacal de La Calera - Frutilla" rows
# (date 2023-10-10, serial 45209) at the top of this market's data block
# (row 489), pushing the existing rows down by 4 and extending the sheet
# from A1:T590 to A1:T594.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 489; this shifts old rows 489-590
# down to 493-594 and keeps the trailing 4 rows (old 587-590) intact at
# the new bottom (591-594).
$ws.Range("A489:A492").EntireRow.Insert()

# Common (constant) values shared by every row in this market/product block
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$fecha = 45209
$codreg = 5
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "`$/bandeja 7 kilos"
$origen = "Provincia de Melipilla"
$kgUnidad = 7

for ($r = 489; $r -le 492; $r++) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

# Row 489: Especial
$ws.Cells.Item(489, 12).Value = "Especial"
$ws.Cells.Item(489, 13).Value = 56
$ws.Cells.Item(489, 14).Value = 11000
$ws.Cells.Item(489, 15).Value = 11000
$ws.Cells.Item(489, 16).Value = 11000
$ws.Cells.Item(489, 19).Value = 1571

# Row 490: Primera
$ws.Cells.Item(490, 12).Value = "Primera"
$ws.Cells.Item(490, 13).Value = 60
$ws.Cells.Item(490, 14).Value = 9000
$ws.Cells.Item(490, 15).Value = 9000
$ws.Cells.Item(490, 16).Value = 9000
$ws.Cells.Item(490, 19).Value = 1286

# Row 491: Segunda
$ws.Cells.Item(491, 12).Value = "Segunda"
$ws.Cells.Item(491, 13).Value = 50
$ws.Cells.Item(491, 14).Value = 7000
$ws.Cells.Item(491, 15).Value = 7000
$ws.Cells.Item(491, 16).Value = 7000
$ws.Cells.Item(491, 19).Value = 1000

# Row 492: Tercera
$ws.Cells.Item(492, 12).Value = "Tercera"
$ws.Cells.Item(492, 13).Value = 40
$ws.Cells.Item(492, 14).Value = 5000
$ws.Cells.Item(492, 15).Value = 5000
$ws.Cells.Item(492, 16).Value = 5000
$ws.Cells.Item(492, 19).Value = 714

# Make sure the sheet's dimension reflects the new extent
$ws.Range("A1:T594").Select()
